$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 date/time values
$ws.Range("D2").Value = "Thu, Sep 28, 2023"
$ws.Range("E2").Value = "9:41 PM"

# Delete rows 3 and 4 (previously had additional attendance records)
$ws.Rows("3:4").Delete()
